$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "Ajo" data block (row 383), pushing
# every existing row below it down by two rows. This mirrors the weekly
# prepend of the latest price observations.
$ws.Rows("383:384").Insert()

# New row 383: $/caja 10 kilos observation for 2023-07-28
$ws.Cells.Item(383,1).Value  = 7
$ws.Cells.Item(383,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(383,3).Value  = "Ñuble"
$ws.Cells.Item(383,4).Value  = (Get-Date -Year 2023 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(383,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(383,5).Value  = 16
$ws.Cells.Item(383,6).Value  = 100112003
$ws.Cells.Item(383,7).Value  = "Ajo"
$ws.Cells.Item(383,8).Value  = "Chino"
$ws.Cells.Item(383,9).Value  = "Primera"
$ws.Cells.Item(383,10).Value = 40
$ws.Cells.Item(383,11).Value = 19000
$ws.Cells.Item(383,12).Value = 19000
$ws.Cells.Item(383,13).Value = 19000
$ws.Cells.Item(383,14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(383,15).Value = "China"
$ws.Cells.Item(383,16).Value = 1900
$ws.Cells.Item(383,17).Value = 10
$ws.Cells.Item(383,18).Value = "Hortaliza"

# New row 384: $/malla 10 kilos observation for 2023-07-28
$ws.Cells.Item(384,1).Value  = 7
$ws.Cells.Item(384,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(384,3).Value  = "Ñuble"
$ws.Cells.Item(384,4).Value  = (Get-Date -Year 2023 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(384,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(384,5).Value  = 16
$ws.Cells.Item(384,6).Value  = 100112003
$ws.Cells.Item(384,7).Value  = "Ajo"
$ws.Cells.Item(384,8).Value  = "Chino"
$ws.Cells.Item(384,9).Value  = "Primera"
$ws.Cells.Item(384,10).Value = 30
$ws.Cells.Item(384,11).Value = 21000
$ws.Cells.Item(384,12).Value = 21000
$ws.Cells.Item(384,13).Value = 21000
$ws.Cells.Item(384,14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(384,15).Value = "China"
$ws.Cells.Item(384,16).Value = 2100
$ws.Cells.Item(384,17).Value = 10
$ws.Cells.Item(384,18).Value = "Hortaliza"
